$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 936.36365
$ws.Range("I39").Value = 951.6667
$ws.Range("J39").Value = 918
$ws.Range("K39").Value = 2855.0001
$ws.Range("L39").Value = 2754
$ws.Range("M39").Value = -2559.0001
$ws.Range("N39").Value = -3346
$ws.Range("H97").Value = 3399
$ws.Range("J97").Value = 3399
$ws.Range("L97").Value = 10197
$ws.Range("N97").Value = -11189
$ws.Range("H111").Value = 1143.6154
$ws.Range("I111").Value = 1077
$ws.Range("J111").Value = 1185.25
$ws.Range("K111").Value = 3231
$ws.Range("L111").Value = 3555.75
$ws.Range("M111").Value = -164
$ws.Range("N111").Value = -9689.75
$ws.Range("H112").Value = 2263.2703
$ws.Range("J112").Value = 2540.0322
$ws.Range("L112").Value = 7620.096600000001
$ws.Range("N112").Value = -9836.096600000001
$ws.Range("H113").Value = 3014.077
$ws.Range("I113").Value = 2398
$ws.Range("J113").Value = 3999.8
$ws.Range("K113").Value = 2398
$ws.Range("L113").Value = 3999.8
$ws.Range("M113").Value = 856
$ws.Range("N113").Value = -10507.8
$ws.Range("H117").Value = 38000
$ws.Range("J117").Value = 38000
$ws.Range("L117").Value = 38000
$ws.Range("N117").Value = -47178
$ws.Range("H132").Value = 18526956
$ws.Range("I132").Value = 27788484
$ws.Range("J132").Value = 3902
$ws.Range("K132").Value = 83365452
$ws.Range("L132").Value = 11706
$ws.Range("M132").Value = -83362922
$ws.Range("N132").Value = -16766
$ws.Range("H135").Value = 1238.7333
$ws.Range("I135").Value = 429.30768
$ws.Range("J135").Value = 6500
$ws.Range("K135").Value = 3863.76912
$ws.Range("L135").Value = 58500
$ws.Range("M135").Value = -1328.76912
$ws.Range("N135").Value = -63570
$ws.Range("H141").Value = 1462.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3274
$ws.Range("I32").Value = 3243.0518
$ws.Range("K32").Value = 3243.0518
$ws.Range("M32").Value = -2956.0518
$ws.Range("H61").Value = 1476.4
$ws.Range("I61").Value = 1109.1428
$ws.Range("J61").Value = 2333.3333
$ws.Range("K61").Value = 1109.1428
$ws.Range("L61").Value = 2333.3333
$ws.Range("M61").Value = -897.1428000000001
$ws.Range("N61").Value = -2757.3333
$ws.Range("H63").Value = 1400
$ws.Range("I63").Value = 1000
$ws.Range("K63").Value = 1000
$ws.Range("M63").Value = -314
$ws.Range("H66").Value = 1400
$ws.Range("I66").Value = 1000
$ws.Range("K66").Value = 5000
$ws.Range("M66").Value = -1568
$ws.Range("H132").Value = 1884.5135
$ws.Range("I132").Value = 1540
$ws.Range("K132").Value = 4620
$ws.Range("M132").Value = -2090
$ws.Range("H136").Value = 1476.4
$ws.Range("I136").Value = 1109.1428
$ws.Range("J136").Value = 2333.3333
$ws.Range("K136").Value = 3327.4284
$ws.Range("L136").Value = 6999.999899999999
$ws.Range("M136").Value = -777.4284000000002
$ws.Range("N136").Value = -12099.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 672.05
$ws.Range("I80").Value = 467.1111
$ws.Range("J80").Value = 839.7273
$ws.Range("K80").Value = 467.1111
$ws.Range("L80").Value = 839.7273
$ws.Range("M80").Value = 530.8888999999999
$ws.Range("N80").Value = -2835.7273
$ws.Range("H83").Value = 672.05
$ws.Range("I83").Value = 467.1111
$ws.Range("J83").Value = 839.7273
$ws.Range("K83").Value = 2335.5555
$ws.Range("L83").Value = 4198.636500000001
$ws.Range("M83").Value = 2656.4445
$ws.Range("N83").Value = -14182.6365
$ws.Range("H134").Value = 7032.75
$ws.Range("I134").Value = 1103.4375
$ws.Range("K134").Value = 3310.3125
$ws.Range("M134").Value = -775.3125
$ws.Range("H135").Value = 33415.5
$ws.Range("J135").Value = 33415.5
$ws.Range("L135").Value = 33415.5
$ws.Range("N135").Value = -43555.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1005
$ws.Range("I12").Value = 1005
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1005
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -835
$ws.Range("H31").Value = 1450.7241
$ws.Range("I31").Value = 1308.28
$ws.Range("J31").Value = 2341
$ws.Range("K31").Value = 1308.28
$ws.Range("L31").Value = 2341
$ws.Range("M31").Value = -1013.28
$ws.Range("N31").Value = -2931
$ws.Range("H34").Value = 1450.7241
$ws.Range("I34").Value = 1308.28
$ws.Range("J34").Value = 2341
$ws.Range("K34").Value = 1308.28
$ws.Range("L34").Value = 2341
$ws.Range("M34").Value = -1106.28
$ws.Range("N34").Value = -2745
$ws.Range("H55").Value = 2000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H132").Value = 2038.3793
$ws.Range("I132").Value = 1463.2106
$ws.Range("J132").Value = 3131.2
$ws.Range("K132").Value = 4389.6318
$ws.Range("L132").Value = 9393.599999999999
$ws.Range("M132").Value = -1859.6318
$ws.Range("N132").Value = -14453.6
$ws.Range("H134").Value = 1246.037
$ws.Range("I134").Value = 1070.5416
$ws.Range("J134").Value = 2650
$ws.Range("K134").Value = 3211.6248
$ws.Range("L134").Value = 7950
$ws.Range("M134").Value = -676.6248000000001
$ws.Range("N134").Value = -13020

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 2714.5
$ws.Range("J106").Value = 2714.5
$ws.Range("L106").Value = 8143.5
$ws.Range("N106").Value = -10035.5
$ws.Range("H121").Value = 1000
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("H136").Value = 1668.8
$ws.Range("J136").Value = 1339.7778
$ws.Range("L136").Value = 4019.3334
$ws.Range("N136").Value = -14219.3334
$ws.Range("H138").Value = 1814.8572
$ws.Range("I138").Value = 1346.1818
$ws.Range("J138").Value = 3533.3333
$ws.Range("K138").Value = 4038.5454
$ws.Range("L138").Value = 10599.9999
$ws.Range("M138").Value = 1101.4546
$ws.Range("N138").Value = -20879.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 500000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H122").Value = 1664.45
$ws.Range("I122").Value = 1900.0714
$ws.Range("J122").Value = 1114.6666
$ws.Range("K122").Value = 5700.2142
$ws.Range("L122").Value = 3343.9998
$ws.Range("M122").Value = -3250.2142
$ws.Range("N122").Value = -8243.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1179.875
$ws.Range("I16").Value = 1179.875
$ws.Range("K16").Value = 1179.875
$ws.Range("M16").Value = -1009.875
$ws.Range("H61").Value = 1796.4445
$ws.Range("I61").Value = 1481.6
$ws.Range("K61").Value = 1481.6
$ws.Range("M61").Value = -1279.6
$ws.Range("H113").Value = 1796.4445
$ws.Range("I113").Value = 1481.6
$ws.Range("K113").Value = 1481.6
$ws.Range("M113").Value = 688.4000000000001
$ws.Range("H132").Value = 27309.281
$ws.Range("I132").Value = 1066.091
$ws.Range("K132").Value = 3198.273
$ws.Range("M132").Value = -668.2729999999997
$ws.Range("H136").Value = 951.0952
$ws.Range("I136").Value = 942.94446
$ws.Range("K136").Value = 2828.83338
$ws.Range("M136").Value = -278.83338

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 14500
$ws.Range("J70").Value = 14500
$ws.Range("L70").Value = 14500
$ws.Range("N70").Value = -15130
$ws.Range("H73").Value = 14500
$ws.Range("J73").Value = 14500
$ws.Range("L73").Value = 14500
$ws.Range("N73").Value = -16684
$ws.Range("H132").Value = 1909.2667
$ws.Range("I132").Value = 1884.1666
$ws.Range("J132").Value = 1946.9166
$ws.Range("K132").Value = 5652.4998
$ws.Range("L132").Value = 5840.7498
$ws.Range("M132").Value = -3122.4998
$ws.Range("N132").Value = -10900.7498
$ws.Range("H136").Value = 531.7692
$ws.Range("I136").Value = 492.75
$ws.Range("K136").Value = 1478.25
$ws.Range("M136").Value = 1071.75
